$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019_PM2.5")

# Add the data row that was removed/emptied from this reduced test fixture.
$ws.Range("A2").Value = 2005
[void]$ws.Range("A2").Select()
